$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" (column G) for the
# 7054fcc4... and b902b015... rows (rows 3 and 4)
$wsOverview.Range("G3").Value = "2016-08-20 22:15:21"
$wsOverview.Range("G4").Value = "2016-08-20 22:15:21"

# zh-cn sheet: Priority (column E) ht -> mt, for rows 3 and 4
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E4").Value = "mt"

# zh-cn sheet: Correspond Handoff Datetime (column H), rows 3 and 4
$wsZhCn.Range("H3").Value = "2016-08-20 22:15:17"
$wsZhCn.Range("H4").Value = "2016-08-20 22:15:17"

# zh-cn sheet: Correspond Handback DateTime (column K), rows 3 and 4
$wsZhCn.Range("K3").Value = "2016-08-20 22:15:31"
$wsZhCn.Range("K4").Value = "2016-08-20 22:15:31"

# de-de sheet: Priority (column E) ht -> mt, for rows 3 and 4
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E4").Value = "mt"

# de-de sheet: Correspond Handoff Datetime (column H), rows 3 and 4
$wsDeDe.Range("H3").Value = "2016-08-20 22:15:21"
$wsDeDe.Range("H4").Value = "2016-08-20 22:15:21"

# de-de sheet: Correspond Handback DateTime (column K), rows 3 and 4
$wsDeDe.Range("K3").Value = "2016-08-20 22:15:37"
$wsDeDe.Range("K4").Value = "2016-08-20 22:15:37"
